# Updates the "HotStock_Top20" sheet data (rows 2-21, columns A-C) to the
# new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("实达集团", "实达集团", "实达集团"),
    @("视觉中国", "视觉中国", "平潭发展"),
    @("榕基软件", "久其软件", "航天发展"),
    @("易点天下", "航天发展", "榕基软件"),
    @("国风新材", "国风新材", "视觉中国"),
    @("航天发展", "易点天下", "久其软件"),
    @("平潭发展", "榕基软件", "国风新材"),
    @("久其软件", "特发信息", "中水渔业"),
    @("赛微电子", "新 华 都", "海南海药"),
    @("蓝色光标", "平潭发展", "大鹏工业"),
    @("特发信息", "龙溪股份", "新华都"),
    @("新 华 都", "北新路桥", "合富中国"),
    @("中水渔业", "赛微电子", "易点天下"),
    @("深中华A", "蓝色光标", "凯美特气"),
    @("中船防务", "贵广网络", "浪潮软件"),
    @("龙溪股份", "中水渔业", "华夏幸福"),
    @("浪潮软件", "欢瑞世纪", "特发信息"),
    @("贵广网络", "深中华A", "华胜天成"),
    @("欢瑞世纪", "郑中设计", "深中华A"),
    @("合富中国", "财信发展", "长城军工")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][2]
}
